# Auto-generated (then hand-reviewed) edit script for cryptos.xlsx
# Applies per-cell text updates matching the target diff, forcing text
# (NumberFormat "@") only for values that Excel would otherwise parse as numbers,
# so genuinely numeric-looking prices stay stored as text (t="inlineStr"-equivalent)
# just like the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextCell {
    param($ws, $ref, $value, $forceText)
    $cell = $ws.Range($ref)
    if ($forceText) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $value
}

Set-TextCell $ws "D2" "26.641.46" $false
Set-TextCell $ws "E2" "  +1.39%  " $false
Set-TextCell $ws "D3" "1.632.52" $false
Set-TextCell $ws "E3" "  +1.61%  " $false
Set-TextCell $ws "E4" "  +0.12%  " $false
Set-TextCell $ws "D5" "212.56" $true
Set-TextCell $ws "E5" "  -0.11%  " $false
Set-TextCell $ws "D6" "0.494" $true
Set-TextCell $ws "E6" "  +1.48%  " $false
Set-TextCell $ws "E7" "  +0.13%  " $false
Set-TextCell $ws "E8" "  +0.88%  " $false
Set-TextCell $ws "E9" "  +1.46%  " $false
Set-TextCell $ws "D10" "19.01" $true
Set-TextCell $ws "E10" "  +3.16%  " $false
Set-TextCell $ws "D11" "0.0841" $true
Set-TextCell $ws "E11" "  +3.42%  " $false
Set-TextCell $ws "D12" "1.862.21" $false
Set-TextCell $ws "E12" "  +1.78%  " $false
Set-TextCell $ws "D13" "1.639.03" $false
Set-TextCell $ws "E13" "  +2.04%  " $false
Set-TextCell $ws "D15" "0.527" $true
Set-TextCell $ws "E15" "  +2.33%  " $false
Set-TextCell $ws "D16" "26.631.70" $false
Set-TextCell $ws "E16" "  +1.50%  " $false
Set-TextCell $ws "D17" "62.91" $true
Set-TextCell $ws "E17" "  +1.36%  " $false
Set-TextCell $ws "D18" "0.0₃0740" $false
Set-TextCell $ws "D19" "208.68" $true
Set-TextCell $ws "E19" "  +4.11%  " $false
Set-TextCell $ws "E20" "  +0.09%  " $false
Set-TextCell $ws "B22" "Avalanche" $false
Set-TextCell $ws "C22" "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax" $false
Set-TextCell $ws "D22" "9.39" $true
Set-TextCell $ws "E22" "  +0.63%  " $false
Set-TextCell $ws "B23" "Chainlink" $false
Set-TextCell $ws "C23" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link" $false
Set-TextCell $ws "D23" "6.17" $true
Set-TextCell $ws "E23" "  +2.61%  " $false
Set-TextCell $ws "E24" "  +2.26%  " $false
Set-TextCell $ws "D25" "146.67" $true
Set-TextCell $ws "E25" "  +1.86%  " $false
Set-TextCell $ws "E26" "  +0.15%  " $false
Set-TextCell $ws "E27" "  -0.92%  " $false
Set-TextCell $ws "D28" "6.77" $true
Set-TextCell $ws "E28" "  +3.04%  " $false
Set-TextCell $ws "D29" "15.35" $true
Set-TextCell $ws "E29" "  +0.93%  " $false
Set-TextCell $ws "E30" "  +5.31%  " $false
Set-TextCell $ws "E31" "  -0.51%  " $false
Set-TextCell $ws "E32" "  +0.94%  " $false
Set-TextCell $ws "D33" "2.95" $true
Set-TextCell $ws "E33" "  -0.06%  " $false
Set-TextCell $ws "E34" "  +0.33%  " $false
Set-TextCell $ws "E35" "  +0.07%  " $false
Set-TextCell $ws "D36" "1.167.71" $false
Set-TextCell $ws "E36" "  +0.19%  " $false
Set-TextCell $ws "E37" "  -1.28%  " $false
Set-TextCell $ws "E38" "  +2.72%  " $false
Set-TextCell $ws "E39" "  +0.13%  " $false
Set-TextCell $ws "B40" "ImmutableX" $false
Set-TextCell $ws "C40" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx" $false
Set-TextCell $ws "D40" "0.503" $true
Set-TextCell $ws "E40" "  +1.01%  " $false
Set-TextCell $ws "B41" "MXToken" $false
Set-TextCell $ws "C41" "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx" $false
Set-TextCell $ws "D41" "2.32" $true
Set-TextCell $ws "E41" "  +0.17%  " $false
Set-TextCell $ws "E42" "  +1.00%  " $false
Set-TextCell $ws "D43" "5.37" $true
Set-TextCell $ws "E43" "  +0.47%  " $false
Set-TextCell $ws "D44" "1.773.03" $false
Set-TextCell $ws "E44" "  +1.87%  " $false
Set-TextCell $ws "D45" "92.01" $true
Set-TextCell $ws "E45" "  -0.15%  " $false
Set-TextCell $ws "D46" "1.54" $true
Set-TextCell $ws "E46" "  +0.19%  " $false
Set-TextCell $ws "E47" "  -2.26%  " $false
Set-TextCell $ws "D48" "54.59" $true
Set-TextCell $ws "E48" "  +0.85%  " $false
Set-TextCell $ws "E49" "  +1.44%  " $false
Set-TextCell $ws "B50" "EnergySwap" $false
Set-TextCell $ws "C50" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens" $false
Set-TextCell $ws "D50" "7.55" $true
Set-TextCell $ws "E50" "  +4.35%  " $false
Set-TextCell $ws "B51" "Mantle" $false
Set-TextCell $ws "C51" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt" $false
Set-TextCell $ws "D51" "0.409" $true
Set-TextCell $ws "E51" "  +0.61%  " $false
